# Workbook and the original "data" worksheet
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Update the "time_taken" column (F2:F18) on the "data" sheet with the
#    refreshed timestamps.
# ---------------------------------------------------------------------------
$data.Range("F2").Value  = "2021-10-05 14:23:16.362917"
$data.Range("F3").Value  = "2021-10-05 14:23:16.362925"
$data.Range("F4").Value  = "2021-10-05 14:23:16.362929"
$data.Range("F5").Value  = "2021-10-05 14:23:16.362932"
$data.Range("F6").Value  = "2021-10-05 14:23:16.362935"
$data.Range("F7").Value  = "2021-10-05 14:23:16.362938"
$data.Range("F8").Value  = "2021-10-05 14:23:16.362941"
$data.Range("F9").Value  = "2021-10-05 14:23:16.362944"
$data.Range("F10").Value = "2021-10-05 14:23:16.362947"
$data.Range("F11").Value = "2021-10-05 14:23:16.362950"
$data.Range("F12").Value = "2021-10-05 14:23:16.362952"
$data.Range("F13").Value = "2021-10-05 14:23:16.362955"
$data.Range("F14").Value = "2021-10-05 14:23:16.362958"
$data.Range("F15").Value = "2021-10-05 14:23:16.362961"
$data.Range("F16").Value = "2021-10-05 14:23:16.362964"
$data.Range("F17").Value = "2021-10-05 14:23:16.362967"
$data.Range("F18").Value = "2021-10-05 14:23:16.362970"

# ---------------------------------------------------------------------------
# 2. Add a new "metadata" worksheet after "data".
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Copy the header cell formatting (bold, bordered, centered) from the
# "data" sheet's header row so the new header reuses the same cell style
# instead of creating a duplicate one.
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Reuse the same style as "data"!A2 (bordered, numeric) for "metadata"!A2.
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)      # xlPasteFormats
$excel.CutCopyMode = $false

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Xeroderma pigmentosum, Trichothiodystrophy or Cockayne syndrome"
$meta.Range("C2").Value = 77
# Force "2.15" to be stored as text, matching the source data which keeps
# the panel version as a string rather than a numeric value.
$meta.Range("D2").Value = "'2.15"
$meta.Range("E2").Value = "2021-09-30T14:47:19.541558Z"
$meta.Range("F2").Value = "2021-10-05 14:23:16.359581"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/77/?format=json"

# Leave the "data" sheet selected/active, matching the original workbook.
$data.Select()
